$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("helpError")

$ws.Range("A7").Value = "Error: Don't know how to add RHS to a theme object"
$ws.Range("B7").Value = "Error: Don't know how to add RHS to an 'x' object"
$ws.Range("C7").Value = "Common in ggplot2 errors: are you trying to add a layer, but forgot a plus sign  at the end of the line somewhere?"

$ws.Range("C8").Select()
